$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 574.5
$ws.Range("I41").Value = 711.1111
$ws.Range("J41").Value = 164.66667
$ws.Range("K41").Value = 711.1111
$ws.Range("L41").Value = 164.66667
$ws.Range("M41").Value = -271.1111
$ws.Range("N41").Value = -1044.66667
$ws.Range("H43").Value = 30262.9
$ws.Range("J43").Value = 18369.666
$ws.Range("L43").Value = 18369.666
$ws.Range("N43").Value = -18507.666
$ws.Range("H59").Value = 1665
$ws.Range("J59").Value = 1665
$ws.Range("L59").Value = 4995
$ws.Range("N59").Value = -6109
$ws.Range("H96").Value = 1084.6666
$ws.Range("I96").Value = 1105.7273
$ws.Range("J96").Value = 1026.75
$ws.Range("K96").Value = 3317.1819
$ws.Range("L96").Value = 3080.25
$ws.Range("M96").Value = -1944.1819
$ws.Range("N96").Value = -5826.25
$ws.Range("H100").Value = 2004.2858
$ws.Range("I100").Value = 2004.2858
$ws.Range("K100").Value = 2004.2858
$ws.Range("M100").Value = -1463.2858
$ws.Range("H107").Value = 4500
$ws.Range("I107").Value = 4000
$ws.Range("K107").Value = 4000
$ws.Range("M107").Value = -2080
$ws.Range("H111").Value = 1043
$ws.Range("I111").Value = 1043
$ws.Range("K111").Value = 3129
$ws.Range("M111").Value = -62
$ws.Range("H125").Value = 1473
$ws.Range("I125").Value = 1234
$ws.Range("K125").Value = 11106
$ws.Range("M125").Value = -8646
$ws.Range("H129").Value = 7799.4
$ws.Range("I129").Value = 7799.4
$ws.Range("K129").Value = 23398.2
$ws.Range("M129").Value = -18398.2
$ws.Range("H131").Value = 24351.2
$ws.Range("J131").Value = 39011
$ws.Range("L131").Value = 117033
$ws.Range("N131").Value = -127113
$ws.Range("H132").Value = 15922162
$ws.Range("I132").Value = 15922162
$ws.Range("K132").Value = 47766486
$ws.Range("M132").Value = -47763956
$ws.Range("H137").Value = 2193.2942
$ws.Range("J137").Value = 1966
$ws.Range("L137").Value = 5898
$ws.Range("N137").Value = -10998
$ws.Range("H138").Value = 3680.6487
$ws.Range("I138").Value = 1539.138
$ws.Range("J138").Value = 5060.7334
$ws.Range("K138").Value = 4617.414
$ws.Range("L138").Value = 15182.2002
$ws.Range("M138").Value = 522.5860000000002
$ws.Range("N138").Value = -25462.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3832.4
$ws.Range("I2").Value = 4097.4
$ws.Range("K2").Value = 4097.4
$ws.Range("M2").Value = -3984.4
$ws.Range("H74").Value = 15817.833
$ws.Range("I74").Value = 2456
$ws.Range("J74").Value = 22498.75
$ws.Range("K74").Value = 2456
$ws.Range("L74").Value = 22498.75
$ws.Range("M74").Value = -1582
$ws.Range("N74").Value = -24246.75
$ws.Range("H77").Value = 15817.833
$ws.Range("I77").Value = 2456
$ws.Range("J77").Value = 22498.75
$ws.Range("K77").Value = 12280
$ws.Range("L77").Value = 112493.75
$ws.Range("M77").Value = -7912
$ws.Range("N77").Value = -121229.75
$ws.Range("H102").Value = 3453.7856
$ws.Range("I102").Value = 3279.4167
$ws.Range("K102").Value = 3279.4167
$ws.Range("M102").Value = -1657.4167
$ws.Range("H116").Value = 3832.4
$ws.Range("I116").Value = 4097.4
$ws.Range("K116").Value = 4097.4
$ws.Range("M116").Value = -1803.4
$ws.Range("H122").Value = 73271.42999999999
$ws.Range("I122").Value = 100799.6
$ws.Range("K122").Value = 302398.8
$ws.Range("M122").Value = -299948.8
$ws.Range("H132").Value = 2721.7058
$ws.Range("I132").Value = 2302.9092
$ws.Range("J132").Value = 3489.5
$ws.Range("K132").Value = 6908.7276
$ws.Range("L132").Value = 10468.5
$ws.Range("M132").Value = -4378.7276
$ws.Range("N132").Value = -15528.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3832.4
$ws.Range("I3").Value = 4097.4
$ws.Range("K3").Value = 4097.4
$ws.Range("M3").Value = -3983.4
$ws.Range("H86").Value = 5441.8667
$ws.Range("I86").Value = 3692.0527
$ws.Range("K86").Value = 3692.0527
$ws.Range("M86").Value = -2569.0527
$ws.Range("H89").Value = 5441.8667
$ws.Range("I89").Value = 3692.0527
$ws.Range("K89").Value = 18460.2635
$ws.Range("M89").Value = -12844.2635
$ws.Range("H94").Value = 1885
$ws.Range("I94").Value = 1332.92
$ws.Range("K94").Value = 1332.92
$ws.Range("M94").Value = -881.9200000000001
$ws.Range("H99").Value = 3942.524
$ws.Range("I99").Value = 1519.5333
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 1519.5333
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -21.53330000000005
$ws.Range("N99").Value = -12996
$ws.Range("H105").Value = 4952.5713
$ws.Range("I105").Value = 5835.6
$ws.Range("K105").Value = 5835.6
$ws.Range("M105").Value = -4088.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 60864.61
$ws.Range("I31").Value = 78660
$ws.Range("J31").Value = 14596.6
$ws.Range("K31").Value = 78660
$ws.Range("L31").Value = 14596.6
$ws.Range("M31").Value = -78365
$ws.Range("N31").Value = -15186.6
$ws.Range("H34").Value = 60864.61
$ws.Range("I34").Value = 78660
$ws.Range("J34").Value = 14596.6
$ws.Range("K34").Value = 78660
$ws.Range("L34").Value = 14596.6
$ws.Range("M34").Value = -78458
$ws.Range("N34").Value = -15000.6
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 2553.7778
$ws.Range("J99").Value = 4004
$ws.Range("K99").Value = 2553.7778
$ws.Range("L99").Value = 4004
$ws.Range("M99").Value = -1055.7778
$ws.Range("N99").Value = -7000
$ws.Range("H105").Value = 2131.5
$ws.Range("I105").Value = 2010.25
$ws.Range("J105").Value = 2252.75
$ws.Range("K105").Value = 2010.25
$ws.Range("L105").Value = 2252.75
$ws.Range("M105").Value = -263.25
$ws.Range("N105").Value = -5746.75
$ws.Range("H107").Value = 1502
$ws.Range("J107").Value = 1226
$ws.Range("L107").Value = 1226
$ws.Range("N107").Value = -5066
$ws.Range("H122").Value = 1088
$ws.Range("I122").Value = 1010
$ws.Range("K122").Value = 3030
$ws.Range("M122").Value = -580
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 2553.7778
$ws.Range("J126").Value = 4004
$ws.Range("K126").Value = 7661.3334
$ws.Range("L126").Value = 12012
$ws.Range("M126").Value = -5191.3334
$ws.Range("N126").Value = -16952
$ws.Range("H134").Value = 8143.0757
$ws.Range("I134").Value = 5774.646
$ws.Range("K134").Value = 17323.938
$ws.Range("M134").Value = -14788.938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 26555.334
$ws.Range("I87").Value = 17799.6
$ws.Range("K87").Value = 53398.8
$ws.Range("M87").Value = -52150.8
$ws.Range("H90").Value = 26555.334
$ws.Range("I90").Value = 17799.6
$ws.Range("K90").Value = 160196.4
$ws.Range("M90").Value = -153956.4
$ws.Range("H122").Value = 2257.8823
$ws.Range("J122").Value = 2083.4167
$ws.Range("L122").Value = 18750.7503
$ws.Range("N122").Value = -23650.7503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3535.6365
$ws.Range("I102").Value = 1699.5714
$ws.Range("J102").Value = 6748.75
$ws.Range("K102").Value = 1699.5714
$ws.Range("L102").Value = 6748.75
$ws.Range("M102").Value = -77.57140000000004
$ws.Range("N102").Value = -9992.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6752.125
$ws.Range("I7").Value = 7227.55
$ws.Range("J7").Value = 4375
$ws.Range("K7").Value = 7227.55
$ws.Range("L7").Value = 4375
$ws.Range("M7").Value = -7115.55
$ws.Range("N7").Value = -4599
$ws.Range("H9").Value = 1137.8
$ws.Range("I9").Value = 172.25
$ws.Range("K9").Value = 172.25
$ws.Range("M9").Value = 51.75
$ws.Range("H35").Value = 3975.7334
$ws.Range("I35").Value = 1738
$ws.Range("K35").Value = 1738
$ws.Range("M35").Value = -1402
$ws.Range("H40").Value = 5582.4287
$ws.Range("I40").Value = 4914.4
$ws.Range("K40").Value = 4914.4
$ws.Range("M40").Value = -4778.4
$ws.Range("H126").Value = 6752.125
$ws.Range("I126").Value = 7227.55
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 21682.65
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -19212.65
$ws.Range("N126").Value = -18065
$ws.Range("H132").Value = 4919.619
$ws.Range("J132").Value = 6736.5
$ws.Range("L132").Value = 20209.5
$ws.Range("N132").Value = -25269.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2347.652
$ws.Range("I107").Value = 2199.8
$ws.Range("K107").Value = 6599.400000000001
$ws.Range("M107").Value = -4679.400000000001
$ws.Range("H122").Value = 2525.276
$ws.Range("I122").Value = 2541.32
$ws.Range("J122").Value = 2425
$ws.Range("K122").Value = 7623.960000000001
$ws.Range("L122").Value = 7275
$ws.Range("M122").Value = -5173.960000000001
$ws.Range("N122").Value = -12175
$ws.Range("H136").Value = 2638.889
$ws.Range("I136").Value = 2495.4348
$ws.Range("J136").Value = 3463.75
$ws.Range("K136").Value = 7486.3044
$ws.Range("L136").Value = 10391.25
$ws.Range("M136").Value = -4936.3044
$ws.Range("N136").Value = -15491.25
